$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map. Values are cryptocurrency prices / hourly volume
# percentages, all stored as plain text in this sheet (Column D = Price,
# Column E = Volume(1h)). Rows 45/46 and 49/50 also swap Coin/Link (B/C).
$updates = @(
    @{Cell = "D2"; Value = "27.723.10"},
    @{Cell = "E2"; Value = "  +0.77%  "},
    @{Cell = "D3"; Value = "1.858.57"},
    @{Cell = "E3"; Value = "  +0.53%  "},
    @{Cell = "E4"; Value = "  -0.88%  "},
    @{Cell = "D5"; Value = "320.06"},
    @{Cell = "E5"; Value = "  -0.20%  "},
    @{Cell = "D6"; Value = "1.018"},
    @{Cell = "E6"; Value = "  -0.73%  "},
    @{Cell = "D7"; Value = "0.4365"},
    @{Cell = "E7"; Value = "  -0.28%  "},
    @{Cell = "D8"; Value = "0.3787"},
    @{Cell = "E8"; Value = "  +0.48%  "},
    @{Cell = "D9"; Value = "0.07420"},
    @{Cell = "E9"; Value = "  +0.31%  "},
    @{Cell = "D10"; Value = "0.8834"},
    @{Cell = "E10"; Value = "  +1.15%  "},
    @{Cell = "D11"; Value = "21.61"},
    @{Cell = "E11"; Value = "  +0.60%  "},
    @{Cell = "D12"; Value = "1.854.69"},
    @{Cell = "E12"; Value = "  +0.16%  "},
    @{Cell = "E13"; Value = "  -0.52%  "},
    @{Cell = "D14"; Value = "6.738"},
    @{Cell = "E14"; Value = "  +0.75%  "},
    @{Cell = "D15"; Value = "0.07134"},
    @{Cell = "E15"; Value = "  -0.42%  "},
    @{Cell = "D16"; Value = "86.84"},
    @{Cell = "E16"; Value = "  +4.81%  "},
    @{Cell = "D17"; Value = "1.023"},
    @{Cell = "E17"; Value = "  -0.83%  "},
    @{Cell = "D18"; Value = "0.000009067"},
    @{Cell = "E18"; Value = "  +0.51%  "},
    @{Cell = "D19"; Value = "1.018"},
    @{Cell = "E19"; Value = "  -0.80%  "},
    @{Cell = "E20"; Value = "  +0.85%  "},
    @{Cell = "D21"; Value = "27.702.60"},
    @{Cell = "E21"; Value = "  +0.62%  "},
    @{Cell = "D22"; Value = "5.290"},
    @{Cell = "E22"; Value = "  +0.79%  "},
    @{Cell = "D23"; Value = "11.15"},
    @{Cell = "E23"; Value = "  -1.44%  "},
    @{Cell = "D24"; Value = "2.102.58"},
    @{Cell = "E24"; Value = "  +1.46%  "},
    @{Cell = "D25"; Value = "2.036"},
    @{Cell = "E25"; Value = "  +6.22%  "},
    @{Cell = "D26"; Value = "156.97"},
    @{Cell = "E26"; Value = "  -0.33%  "},
    @{Cell = "D27"; Value = "18.71"},
    @{Cell = "E27"; Value = "  -0.07%  "},
    @{Cell = "D28"; Value = "1.999"},
    @{Cell = "E28"; Value = "  +1.43%  "},
    @{Cell = "D29"; Value = "5.356"},
    @{Cell = "E29"; Value = "  +1.80%  "},
    @{Cell = "D30"; Value = "120.42"},
    @{Cell = "E30"; Value = "  +2.85%  "},
    @{Cell = "D31"; Value = "0.09054"},
    @{Cell = "E31"; Value = "  +0.14%  "},
    @{Cell = "D32"; Value = "1.220"},
    @{Cell = "E32"; Value = "  +1.99%  "},
    @{Cell = "D33"; Value = "0.7704"},
    @{Cell = "E33"; Value = "  +1.23%  "},
    @{Cell = "D34"; Value = "3.031"},
    @{Cell = "E34"; Value = "  +5.41%  "},
    @{Cell = "E35"; Value = "  +0.85%  "},
    @{Cell = "E36"; Value = "  -0.73%  "},
    @{Cell = "D38"; Value = "0.01976"},
    @{Cell = "E38"; Value = "  +0.17%  "},
    @{Cell = "E39"; Value = "  +0.10%  "},
    @{Cell = "D40"; Value = "2.869"},
    @{Cell = "E40"; Value = "  +2.67%  "},
    @{Cell = "D41"; Value = "0.5190"},
    @{Cell = "E41"; Value = "  +0.93%  "},
    @{Cell = "D42"; Value = "6.950"},
    @{Cell = "E42"; Value = "  +3.21%  "},
    @{Cell = "D43"; Value = "0.1678"},
    @{Cell = "E43"; Value = "  +0.32%  "},
    @{Cell = "D44"; Value = "8.687"},
    @{Cell = "E44"; Value = "  +2.51%  "},
    @{Cell = "B45"; Value = "Quant"},
    @{Cell = "C45"; Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"},
    @{Cell = "D45"; Value = "110.06"},
    @{Cell = "E45"; Value = "  +1.34%  "},
    @{Cell = "B46"; Value = "EnergySwap"},
    @{Cell = "C46"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"},
    @{Cell = "D46"; Value = "10.72"},
    @{Cell = "E46"; Value = "  +1.39%  "},
    @{Cell = "D47"; Value = "1.711"},
    @{Cell = "E47"; Value = "  +0.29%  "},
    @{Cell = "D48"; Value = "0.4724"},
    @{Cell = "E48"; Value = "  +1.75%  "},
    @{Cell = "B49"; Value = "PaxDollar"},
    @{Cell = "C49"; Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"},
    @{Cell = "D49"; Value = "1.020"},
    @{Cell = "E49"; Value = "  -0.86%  "},
    @{Cell = "B50"; Value = "Cronos"},
    @{Cell = "C50"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"},
    @{Cell = "D50"; Value = "0.06500"},
    @{Cell = "E50"; Value = "  +1.67%  "},
    @{Cell = "D51"; Value = "1.854"},
    @{Cell = "E51"; Value = "  +0.42%  "}
)

foreach ($u in $updates) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # instead of auto-converting numeric-looking strings (e.g. "320.06")
    # into a Number cell.
    $ws.Range($u.Cell).Value = "'" + $u.Value
    # Clear the resulting quote-prefix style so the cell keeps its original
    # (default/General) formatting, matching the source data.
    $ws.Range($u.Cell).Style = "Normal"
}
